$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "List name" column (A) is being removed entirely; every other column
# shifts left by one (B->A, C->B, D->C). Deleting the whole column re-uses
# the existing cell contents/styles of the old B/C/D columns, which already
# line up with the target layout.
$ws.Columns.Item(1).Delete()

# Old column B ("Task title") becomes the new column A header, but its text
# changes to "Title".
$ws.Range("A1").Value = "Title"

# New column A width (was two columns: 21.0 and 42.29 wide) becomes a
# single ~34.43-character-wide column.
$ws.Columns.Item(1).ColumnWidth = 33.6667
